# All Country Files Saved And Formatted
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header rename: "RowNo." -> "MasterSheet RowNo."
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 12).Value = "MasterSheet RowNo."

# ---------------------------------------------------------------------------
# 2. Fill in the previously-missing TotalConfirmedNewCases (G) /
#    TotalNewDeaths (I) columns for the existing rows, and append the new
#    row 9 of data.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 7).Value  = 1
$ws.Cells.Item(2, 9).Value  = 0

$ws.Cells.Item(3, 7).Value  = 0
$ws.Cells.Item(3, 9).Value  = 0

$ws.Cells.Item(4, 7).Value  = 1
$ws.Cells.Item(4, 9).Value  = 0

$ws.Cells.Item(5, 7).Value  = 0
$ws.Cells.Item(5, 9).Value  = 0

$ws.Cells.Item(6, 7).Value  = 0
$ws.Cells.Item(6, 9).Value  = 0

$ws.Cells.Item(7, 7).Value  = 0
$ws.Cells.Item(7, 9).Value  = 0

$ws.Cells.Item(8, 7).Value  = 0
$ws.Cells.Item(8, 9).Value  = 0

# New row 9
$ws.Cells.Item(9, 1).Value  = 71
$ws.Cells.Item(9, 2).Value  = 21
$ws.Cells.Item(9, 3).Value  = "LATIN AMER. & CARIB    "
$ws.Cells.Item(9, 4).Value  = 43921
$ws.Cells.Item(9, 5).Value  = "Belize"
$ws.Cells.Item(9, 6).Value  = 3
$ws.Cells.Item(9, 7).Value  = 1
$ws.Cells.Item(9, 8).Value  = 0
$ws.Cells.Item(9, 9).Value  = 0
$ws.Cells.Item(9, 10).Value = "Local transmission"
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 5367

# ---------------------------------------------------------------------------
# 3. Column widths: A:O all 27 "characters" wide. Excel's ColumnWidth
#    property is expressed in a slightly different unit than the raw XML
#    width (an offset of 5/6 character for the cell margin), so back it out
#    here so the persisted <col width="..."/> comes out to an exact 27.
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").EntireColumn.ColumnWidth = 26.1666666666667

# ---------------------------------------------------------------------------
# 4. Styling. Two new cell formats are introduced:
#      - center/center alignment (applies to every used cell)
#      - center/center alignment + a "yyyy-mm-dd;" date format (Date column)
#    Build each format once on a helper cell that's already going to be part
#    of the final used range (so no stray styled cell is left outside the
#    sheet's dimension), then broadcast it with copy/paste-special so every
#    target cell lands on the SAME style index instead of each property
#    assignment minting its own intermediate, unused entry in cellXfs.
# ---------------------------------------------------------------------------

# 4a. plain center/center format, built on O9 (itself one of the new filler
#     cells that should end up with this exact style)
$plain = $ws.Cells.Item(9, 15)
$plain.HorizontalAlignment = -4108
$plain.VerticalAlignment = -4108
$plain.Copy()
$ws.Range("A1:O9").PasteSpecial(-4122)

# 4b. center/center + date format, built on N9 (also one of the filler
#     cells, temporarily reused as the template and then restored below)
$dateTemplate = $ws.Cells.Item(9, 14)
$dateTemplate.NumberFormat = "yyyy-mm-dd;"
$dateTemplate.HorizontalAlignment = -4108
$dateTemplate.VerticalAlignment = -4108
$dateTemplate.Copy()
$ws.Range("D1:D9").PasteSpecial(-4122)

# restore N9 itself back to the plain format (it's just an empty filler
# cell, not part of the Date column)
$plain.Copy()
$dateTemplate.PasteSpecial(-4122)

Write-Host "All Country Files Saved And Formatted"
